$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B65").Value = "Emanuel y Jorgelina"
$ws.Range("C65").Value = 1176381379
$ws.Range("D65").Value = "115 nº 1140"

$ws.Range("D65").Select() | Out-Null
